# Apply changes to Donnees_Projet_Optimisation.xlsx:
#  - Remove the leading blank spacer row above the header on the data sheet
#    (this shifts the header and all data rows up by one row)
#  - Update the selected cell to R1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tableaux2et3et4")

# Before edit, layout is:
#   Row 1 : blank spacer row
#   Row 2 : header
#   Row 3..14 : data rows (Unit 1 .. Unit 12)

# Delete the leading blank spacer row (row 1); this shifts everything up by
# one row, so the header becomes row 1 and the data rows become rows 2..13.
$ws.Rows.Item(1).Delete() | Out-Null

# Update the selection shown in the workbook to R1
$ws.Range("R1").Select() | Out-Null
